# Update countries & provincias Spain
# Applies the 18-Jul-2020 01:11 -> 02:28 data refresh to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp header (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 02:28"

# --- Straightforward numeric refreshes (country keeps its row) ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3766593
$ws.Range("C4").Value = 71568
$ws.Range("D4").Value = 1733847
$ws.Range("E4").Value = 1890770
$ws.Range("G4").Value = 858
$ws.Range("H4").Value = 141976

# Row 19: Alemania
$ws.Range("B19").Value = 202345
$ws.Range("C19").Value = 509
$ws.Range("E19").Value = 6285

# Row 24: Canada
$ws.Range("B24").Value = 109669
$ws.Range("C24").Value = 405
$ws.Range("D24").Value = 96689
$ws.Range("E24").Value = 4141

# Row 42: Panama
$ws.Range("B42").Value = 51408
$ws.Range("C42").Value = 1035
$ws.Range("D42").Value = 26520
$ws.Range("E42").Value = 23850
$ws.Range("G42").Value = 38
$ws.Range("H42").Value = 1038

# Row 167: Guyana
$ws.Range("B167").Value = 320
$ws.Range("C167").Value = 5
$ws.Range("E167").Value = 145

# --- Rank swaps: a country's update pushes it above its former neighbour ---

# Rows 68/69: Chequia overtakes Costa de Marfil
$ws.Range("A68").Value = "Chequia"
$ws.Range("B68").Value = 13742
$ws.Range("C68").Value = 130
$ws.Range("D68").Value = 8725
$ws.Range("E68").Value = 4659
$ws.Range("G68").Value = 3
$ws.Range("H68").Value = 358

$ws.Range("A69").Value = "Costa de Marfil"
$ws.Range("B69").Value = 13696
$ws.Range("C69").Value = 142
$ws.Range("D69").Value = 7607
$ws.Range("E69").Value = 6002
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 87

# Rows 144/145: Uruguay overtakes Republica de Chipre
$ws.Range("A144").Value = "Uruguay"
$ws.Range("B144").Value = 1037
$ws.Range("C144").Value = 11
$ws.Range("D144").Value = 917
$ws.Range("E144").Value = 88
$ws.Range("H144").Value = 32

$ws.Range("A145").Value = "Republica de Chipre"
$ws.Range("B145").Value = 1033
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 845
$ws.Range("E145").Value = 169
$ws.Range("H145").Value = 19

# Rows 190/191: Antigua y Barbuda overtakes Islas Turcas y Caicos
$ws.Range("A190").Value = "Antigua y Barbuda"
$ws.Range("B190").Value = 76
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 57
$ws.Range("E190").Value = 16
$ws.Range("H190").Value = 3

$ws.Range("A191").Value = "Islas Turcas y Caicos"
$ws.Range("B191").Value = 74
$ws.Range("C191").Value = 2
$ws.Range("D191").Value = 12
$ws.Range("E191").Value = 60
$ws.Range("H191").Value = 2

# Rows 193/194: San Martin (Parte Francesa) overtakes Macao
$ws.Range("A193").Value = "San Martin (Parte Francesa)"
$ws.Range("B193").Value = 46
$ws.Range("C193").Value = 3
$ws.Range("D193").Value = 39
$ws.Range("E193").Value = 4
$ws.Range("H193").Value = 3

$ws.Range("A194").Value = "Macao"
$ws.Range("B194").Value = 46
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 45
$ws.Range("E194").Value = 1
$ws.Range("H194").Value = 0
